# Apply updated cryptocurrency price/volume figures to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.405.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.024.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.30%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.020.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.56%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.433"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.524.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.634.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.030.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "441.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.667"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0956"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("E34").Value = "  -2.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.965"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0695"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0369"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("E41").Value = "  -2.14%  "
$ws.Range("E42").Value = "  -7.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.666.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "122.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  -4.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.78%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.83%  "
